# Insert a new data row before the current row 315 ("Vega Central Mapocho
# de Santiago" / Ciboulette weekly series), pushing the existing rows
# 315-335 down to 316-336, and fill in the new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 315:335 down to 316:336 by inserting a blank row at 315.
$ws.Rows(315).Insert()

# Populate the newly inserted row 315 with the new weekly record.
$ws.Cells.Item(315, 1).Value = 9
$ws.Cells.Item(315, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(315, 3).Value = "Metropolitana"
$ws.Cells.Item(315, 4).Value = 44585
$ws.Cells.Item(315, 5).Value = 13
$ws.Cells.Item(315, 6).Value = 100112039
$ws.Cells.Item(315, 7).Value = "Ciboulette"
$ws.Cells.Item(315, 8).Value = "Sin especificar"
$ws.Cells.Item(315, 9).Value = "Primera"
$ws.Cells.Item(315, 10).Value = 106
$ws.Cells.Item(315, 11).Value = 1000
$ws.Cells.Item(315, 12).Value = 1200
$ws.Cells.Item(315, 13).Value = 1100
$ws.Cells.Item(315, 14).Value = "`$/docena de atados"
$ws.Cells.Item(315, 15).Value = "Región Metropolitana"
$ws.Cells.Item(315, 16).Value = 367
$ws.Cells.Item(315, 17).Value = 3
$ws.Cells.Item(315, 18).Value = "Hortaliza"

# Match the date number format already used by the other rows in column D.
$ws.Cells.Item(315, 4).NumberFormat = $ws.Cells.Item(316, 4).NumberFormat
